# Daily attendance processing - 2025-11-10 07:23:11
# Reorder the "Recorded By" text in column G from "dnasr281@gmail.com, System"
# to "System, dnasr281@gmail.com" for every row where it matches exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G ("Recorded By")
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
